$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "reason"

$ws.Range("B2").Value = 100
$ws.Range("C2").Value = "`nThe job of NLP engineer remote is highly suitable for you as it requires the skills of NLP, Pytorch, Computer Vision, and Python, all of which you possess. The job also involves analyzing and preprocessing large-scale text data and developing prototypes, which aligns with your experience in LLMGuard and Literature Society IITJ Website projects. The high score of 100 indicates a strong match between your skills and the job requirements."

$ws.Range("B3").Value = 75
$ws.Range("C3").Value = "The job requires skills in MongoDB, ReactJS, JavaScript, Web Development, and NodeJS, which align with your experience and projects. The job offers the opportunity to work on high-quality web/mobile applications and suggests new features for product enhancement, which aligns with your skills in ReactJS and JavaScript. The job has a moderate score as it is a good fit for you, but there may be other candidates with more relevant experience."

$ws.Range("B4").Value = 62
$ws.Range("C4").Value = "`nThe job has a moderate score of 62, indicating that it is neither highly suitable nor unsuitable for you. This is because the job requires skills in ReactJS, JavaScript, CSS, Frontend Development, and NextJS, which you have experience in through your projects. However, the job is an unpaid internship and for only 3 months, which may not be a good fit for your long-term career goals. "
